$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D2:D51) to Text format so values like "1.000",
# "23.00" and two-dot big numbers like "28.116.96" are preserved exactly
# as strings instead of being auto-coerced to numbers by Excel.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

# --- Coin name (B) changes ---
$ws.Range("B36").Value = "VeChain"
$ws.Range("B37").Value = "Hedera"
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("B49").Value = "Quant"

# --- Link (C) changes ---
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"

# --- Price (D) changes ---
$ws.Range("D2").Value = "28.116.96"
$ws.Range("D3").Value = "1.819.31"
$ws.Range("D4").Value = "1.001"
$ws.Range("D5").Value = "337.83"
$ws.Range("D7").Value = "0.4224"
$ws.Range("D8").Value = "0.3527"
$ws.Range("D9").Value = "45.56"
$ws.Range("D10").Value = "1.164"
$ws.Range("D11").Value = "0.07547"
$ws.Range("D12").Value = "23.00"
$ws.Range("D13").Value = "0.9998"
$ws.Range("D14").Value = "6.318"
$ws.Range("D15").Value = "7.313"
$ws.Range("D16").Value = "1.811.31"
$ws.Range("D17").Value = "0.00001096"
$ws.Range("D18").Value = "0.06702"
$ws.Range("D19").Value = "82.97"
$ws.Range("D20").Value = "1.000"
$ws.Range("D21").Value = "17.48"
$ws.Range("D22").Value = "6.408"
$ws.Range("D23").Value = "28.120.94"
$ws.Range("D24").Value = "11.95"
$ws.Range("D25").Value = "2.403"
$ws.Range("D26").Value = "2.517"
$ws.Range("D27").Value = "20.89"
$ws.Range("D28").Value = "156.56"
$ws.Range("D29").Value = "2.019.81"
$ws.Range("D30").Value = "1.321"
$ws.Range("D31").Value = "133.62"
$ws.Range("D32").Value = "4.083"
$ws.Range("D33").Value = "6.042"
$ws.Range("D34").Value = "0.09163"
$ws.Range("D35").Value = "12.45"
$ws.Range("D36").Value = "0.02361"
$ws.Range("D37").Value = "0.06367"
$ws.Range("D38").Value = "0.6708"
$ws.Range("D39").Value = "5.268"
$ws.Range("D40").Value = "0.2174"
$ws.Range("D42").Value = "1.221"
$ws.Range("D43").Value = "8.184"
$ws.Range("D44").Value = "14.30"
$ws.Range("D46").Value = "0.6193"
$ws.Range("D47").Value = "3.876"
$ws.Range("D48").Value = "2.071"
$ws.Range("D49").Value = "128.70"
$ws.Range("D50").Value = "1.190"
$ws.Range("D51").Value = "0.07129"

# --- Volume(1h) (E) changes ---
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("E3").Value = "  +2.72%  "
$ws.Range("E4").Value = "  -0.82%  "
$ws.Range("E5").Value = "  -0.17%  "
$ws.Range("E6").Value = "  -0.49%  "
$ws.Range("E7").Value = "  +10.37%  "
$ws.Range("E8").Value = "  +3.21%  "
$ws.Range("E9").Value = "  -2.78%  "
$ws.Range("E10").Value = "  +1.35%  "
$ws.Range("E11").Value = "  +2.17%  "
$ws.Range("E12").Value = "  -2.23%  "
$ws.Range("E13").Value = "  -0.69%  "
$ws.Range("E14").Value = "  -1.59%  "
$ws.Range("E15").Value = "  +0.11%  "
$ws.Range("E16").Value = "  +1.74%  "
$ws.Range("E17").Value = "  +1.82%  "
$ws.Range("E18").Value = "  +0.38%  "
$ws.Range("E19").Value = "  +0.82%  "
$ws.Range("E20").Value = "  -0.48%  "
$ws.Range("E21").Value = "  +0.53%  "
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("E23").Value = "  +0.19%  "
$ws.Range("E24").Value = "  -1.01%  "
$ws.Range("E25").Value = "  +0.53%  "
$ws.Range("E26").Value = "  +4.81%  "
$ws.Range("E27").Value = "  +0.90%  "
$ws.Range("E28").Value = "  +1.77%  "
$ws.Range("E29").Value = "  +1.99%  "
$ws.Range("E30").Value = "  -6.36%  "
$ws.Range("E31").Value = "  -1.07%  "
$ws.Range("E32").Value = "  +1.33%  "
$ws.Range("E33").Value = "  -0.36%  "
$ws.Range("E34").Value = "  +2.90%  "
$ws.Range("E35").Value = "  -2.15%  "
$ws.Range("E36").Value = "  -1.47%  "
$ws.Range("E37").Value = "  +0.34%  "
$ws.Range("E38").Value = "  -1.47%  "
$ws.Range("E39").Value = "  -1.09%  "
$ws.Range("E40").Value = "  +0.82%  "
$ws.Range("E41").Value = "  +0.57%  "
$ws.Range("E42").Value = "  -1.75%  "
$ws.Range("E43").Value = "  -0.94%  "
$ws.Range("E44").Value = "  +1.22%  "
$ws.Range("E45").Value = "  -0.46%  "
$ws.Range("E46").Value = "  -0.94%  "
$ws.Range("E47").Value = "  +0.49%  "
$ws.Range("E48").Value = "  +0.37%  "
$ws.Range("E49").Value = "  -2.96%  "
$ws.Range("E50").Value = "  +0.07%  "
$ws.Range("E51").Value = "  -5.10%  "

# Restore the original (default) cell style on the Price column now that
# the text values are set, so no stray number-format style lingers on the
# cells themselves (mirrors the untouched D6 cell's style).
$priceRange.Style = $ws.Range("D6").Style
